$d = $word.ActiveDocument

function New-ListParagraph($level) {
    $base = $d.Paragraphs($d.Paragraphs.Count).Range
    $base.Collapse(0)
    $base.InsertParagraphAfter()
    $p = $d.Paragraphs($d.Paragraphs.Count)
    $p.Range.ListFormat.ListLevelNumber = $level
    return $p
}

function Add-RunText($para, $text) {
    $r = $para.Range
    $r.Collapse(0)
    $r.InsertBefore($text)
}

# --- New paragraph 0: ilvl=0 numId=4 ---
$para = New-ListParagraph 1

# --- New paragraph 1: ilvl=0 numId=4 ---
$para = New-ListParagraph 1
Add-RunText $para "Quarterly visitor nights in Australia data"

# --- New paragraph 2: ilvl=1 numId=4 ---
$para = New-ListParagraph 2
Add-RunText $para "This time series has quarterly frequency. The patterns in the data suggest strong seasonality and an upward trend. The seasonal variation is much "
Add-RunText $para "larger than the rise in the trend."
Add-RunText $para " Also, it seems that the data requires some transformation to stabilize the variance."

# --- New paragraph 3: ilvl=1 numId=4 ---
$para = New-ListParagraph 2
Add-RunText $para "The correlation between present values and past values is statistically significant for the seasonal lags (1-3) and for adjacent lags to seasonal components."

# --- New paragraph 4: ilvl=1 numId=4 ---
$para = New-ListParagraph 2
Add-RunText $para "The PACF suggests significant lags 1 and 2 and significant seasonal lags 1 and 2 as well."

# --- New paragraph 5: ilvl=1 numId=4 ---
$para = New-ListParagraph 2
Add-RunText $para "A "
Add-RunText $para "Box-Cox transformation and one seasonal differencing"
Add-RunText $para " seem to "
Add-RunText $para "create a stationary time-series. We fail to reject the null hypothesis of the KPSS unit root test."
Add-RunText $para " A good model seems to be SARIMA (1,0,1), (1,1,1)"

# --- New paragraph 6: ilvl=1 numId=4 ---
$para = New-ListParagraph 2
Add-RunText $para "No, it chose an SARIMA (1,0,0) (0,1,1). An ARIMA model with one non-seasonal autoregressive component, but no moving average or differencing; and a seasonal part comprised of no autoregressive terms, but that differences the data once and also contains one moving average component. Using cross-validation, the model chosen by auto-arima has a smaller MSE."

# --- New paragraph 7: ilvl=1 numId=4 ---
$para = New-ListParagraph 2
Add-RunText $para "Done in notebook"

# --- New paragraph 8: ilvl=0 numId=4 ---
$para = New-ListParagraph 1
Add-RunText $para "Usmelec series"

# --- New paragraph 9: ilvl=1 numId=4 ---
$para = New-ListParagraph 2
Add-RunText $para "The data clearly has an upward trend."

# --- New paragraph 10: ilvl=1 numId=4 ---
$para = New-ListParagraph 2
Add-RunText $para "The plot of the data shows that the seasonal variation increases with time. A Box-Cox with lambda -0.57 will be used."

# --- New paragraph 11: ilvl=1 numId=4 ---
$para = New-ListParagraph 2
Add-RunText $para "Using the ndiff command after taking first-seasonal differencing seems to suffice to create a stationary series."

# --- New paragraph 12: ilvl=1 numId=4 ---
$para = New-ListParagraph 2
Add-RunText $para "Using auto-arima, we found an ARIMA (1,1,1)"
Add-RunText $para " (2,1,1) model, which has a lower AICc than the model fitted manually ARIMA(1,1,3)(0,1,3)"

# --- New paragraph 13: ilvl=1 numId=4 ---
$para = New-ListParagraph 2
Add-RunText $para "The residuals resemble white-noise."

# --- New paragraph 14: ilvl=1 numId=4 ---
$para = New-ListParagraph 2
Add-RunText $para "Checking against actual values, the model "
Add-RunText $para "has the following error metrics"

# --- New paragraph 15: ilvl=2 numId=4 ---
$para = New-ListParagraph 3
Add-RunText $para "MAE: 8.17"

# --- New paragraph 16: ilvl=2 numId=4 ---
$para = New-ListParagraph 3
Add-RunText $para "MAPE: 2.41"

# --- New paragraph 17: ilvl=2 numId=4 ---
$para = New-ListParagraph 3
Add-RunText $para "RMSE: 10.33"

# --- New paragraph 18: ilvl=1 numId=4 ---
$para = New-ListParagraph 2
Add-RunText $para "Around 2020."

# Re-create the _GoBack bookmark (collapsed) at the very end of the document.
$endRng = $d.Paragraphs($d.Paragraphs.Count).Range
$endRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRng)
